$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("prep_sheet")
$ws.Activate()

# New product catalog value for the "12pm HEARTY SALADS" qty column
# (rows 24-26 switch from the shared "1 quart" string to a new
# "1 pint, raw" string).
$ws.Range("E24").Value = "1 pint, raw"
$ws.Range("E25").Value = "1 pint, raw"
$ws.Range("E26").Value = "1 pint, raw"

# Row 30 now needs extra height to fit the updated text.
$ws.Rows.Item(30).RowHeight = 24

# Reset the view: scroll back to the top of the sheet and move the
# active selection to D38.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D38").Select()
